$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3 and 4 (A, B/C text, D values change)
$ws.Cells.Item(3, 1).Value = 12
$ws.Cells.Item(3, 2).Value = "test3"
$ws.Cells.Item(3, 3).Value = "ijkl"
$ws.Cells.Item(3, 4).Value = 1

$ws.Cells.Item(4, 1).Value = 14
$ws.Cells.Item(4, 2).Value = "test2"
$ws.Cells.Item(4, 3).Value = "efgh"
$ws.Cells.Item(4, 4).Value = 2

# Add new rows 5 and 6
$ws.Cells.Item(5, 1).Value = 15
$ws.Cells.Item(5, 2).Value = "test1"
$ws.Cells.Item(5, 3).Value = "dsf"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 1

$ws.Cells.Item(6, 1).Value = 16
$ws.Cells.Item(6, 2).Value = "test2"
$ws.Cells.Item(6, 3).Value = "sefse"
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 1
